$wb = $excel.ActiveWorkbook

# The values go on the "CNN" sheet (first sheet) which holds the
# "Single Channel (CNN)" results in column B.
$ws = $wb.Worksheets.Item("CNN")
$ws.Activate()

$ws.Range("B3").Value = 0.712
$ws.Range("B4").Value = 0.6978
$ws.Range("B5").Value = 0.8192
$ws.Range("B6").Value = 0.8481

# Update the active selection to match the recorded state after editing.
$ws.Range("B11").Select()
